$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update metadata rows 2-4 for columns C, E, F to reflect the newly curated dimensions
$ws.Range("C2").Value = "sdmx-dimension:refArea"
$ws.Range("E2").Value = "iaest-measure:sexo"
$ws.Range("F2").Value = "iaest-measure:jefe-explotacion"

$ws.Range("C3").Value = "dim"
$ws.Range("E3").Value = "medida"
$ws.Range("F3").Value = "medida"

$ws.Range("C4").Value = "URI-Municipio"
$ws.Range("E4").Value = "xsd:int"
$ws.Range("F4").Value = "xsd:int"

# Remove the no-longer-needed mapping file references in row 5
$ws.Range("E5").Clear()
$ws.Range("F5").Clear()
